$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the FDI row (row 12): series_id, series_name, and source_id change
# because the data now comes from the UNCTAD source instead of the
# U.S. Direct Investment Abroad source.
$ws.Range("A12").Value = 1665
$ws.Range("B12").Value = 'Foreign direct investment (FDI)  inflows ($US millions)'
$ws.Range("D12").Value = 32

# Update the active selection to D13 (matches the saved selection in the file)
$ws.Range("D13").Select()

# Update the Excel application window position / size (reflected in bookViews)
$aw = $excel.ActiveWindow
$aw.Left = 1170
$aw.Top = 960
$aw.Width = 23700
$aw.Height = 15240

$wb.Save()
